$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.837.38"
Set-TextValue $ws.Range("E2") "  +2.06%  "

Set-TextValue $ws.Range("D3") "3.787.94"
Set-TextValue $ws.Range("E3") "  +0.88%  "

Set-TextValue $ws.Range("E4") "  -0.10%  "

Set-TextValue $ws.Range("D5") "597.37"
Set-TextValue $ws.Range("E5") "  +0.41%  "

Set-TextValue $ws.Range("D6") "170.25"
Set-TextValue $ws.Range("E6") "  +0.29%  "

Set-TextValue $ws.Range("D7") "3.784.85"
Set-TextValue $ws.Range("E7") "  +0.87%  "

Set-TextValue $ws.Range("E8") "  -0.07%  "

Set-TextValue $ws.Range("E9") "  -0.47%  "

Set-TextValue $ws.Range("E10") "  -1.47%  "

Set-TextValue $ws.Range("E11") "  +0.98%  "

Set-TextValue $ws.Range("D13") "0.0000263"
Set-TextValue $ws.Range("E13") "  -3.51%  "

Set-TextValue $ws.Range("D14") "36.90"
Set-TextValue $ws.Range("E14") "  +0.58%  "

Set-TextValue $ws.Range("D15") "4.422.34"
Set-TextValue $ws.Range("E15") "  +0.82%  "

Set-TextValue $ws.Range("D16") "3.788.70"
Set-TextValue $ws.Range("E16") "  +0.64%  "

Set-TextValue $ws.Range("D17") "68.817.14"
Set-TextValue $ws.Range("E17") "  +1.90%  "

Set-TextValue $ws.Range("D18") "18.26"
Set-TextValue $ws.Range("E18") "  -3.26%  "

Set-TextValue $ws.Range("D19") "7.08"
Set-TextValue $ws.Range("E19") "  -2.16%  "

Set-TextValue $ws.Range("E20") "  -0.09%  "

Set-TextValue $ws.Range("D21") "11.07"
Set-TextValue $ws.Range("E21") "  +4.94%  "

Set-TextValue $ws.Range("D22") "469.93"
Set-TextValue $ws.Range("E22") "  +0.21%  "

Set-TextValue $ws.Range("D23") "0.707"
Set-TextValue $ws.Range("E23") "  -1.79%  "

Set-TextValue $ws.Range("D24") "84.90"
Set-TextValue $ws.Range("E24") "  +1.31%  "

Set-TextValue $ws.Range("E25") "  -3.22%  "

Set-TextValue $ws.Range("E26") "  +0.51%  "

Set-TextValue $ws.Range("D27") "12.24"
Set-TextValue $ws.Range("E27") "  +0.54%  "

Set-TextValue $ws.Range("D28") "10.21"
Set-TextValue $ws.Range("E28") "  -1.59%  "

Set-TextValue $ws.Range("D30") "3.934.71"
Set-TextValue $ws.Range("E30") "  +0.67%  "

Set-TextValue $ws.Range("D31") "2.82"
Set-TextValue $ws.Range("E31") "  -3.15%  "

Set-TextValue $ws.Range("D32") "7.44"
Set-TextValue $ws.Range("E32") "  -3.08%  "

Set-TextValue $ws.Range("D33") "2.25"
Set-TextValue $ws.Range("E33") "  +0.04%  "

Set-TextValue $ws.Range("D34") "30.23"
Set-TextValue $ws.Range("E34") "  -0.73%  "

Set-TextValue $ws.Range("D35") "9.38"
Set-TextValue $ws.Range("E35") "  +2.78%  "

Set-TextValue $ws.Range("D37") "3.741.92"
Set-TextValue $ws.Range("E37") "  +0.56%  "

Set-TextValue $ws.Range("E38") "  -3.16%  "

Set-TextValue $ws.Range("E39") "  -10.22%  "

Set-TextValue $ws.Range("E40") "  +1.57%  "

Set-TextValue $ws.Range("E41") "  +0.78%  "

Set-TextValue $ws.Range("D42") "5.84"
Set-TextValue $ws.Range("E42") "  -0.78%  "

Set-TextValue $ws.Range("E43") "  -0.08%  "

Set-TextValue $ws.Range("D44") "0.310"

Set-TextValue $ws.Range("E45") "  +0.00%  "

Set-TextValue $ws.Range("D46") "1.98"
Set-TextValue $ws.Range("E46") "  +1.07%  "

Set-TextValue $ws.Range("D47") "43.98"
Set-TextValue $ws.Range("E47") "  +13.20%  "

Set-TextValue $ws.Range("E48") "  -1.40%  "

Set-TextValue $ws.Range("D49") "46.16"
Set-TextValue $ws.Range("E49") "  +0.73%  "

Set-TextValue $ws.Range("D50") "400.88"
Set-TextValue $ws.Range("E50") "  +0.81%  "

Set-TextValue $ws.Range("D51") "145.65"
Set-TextValue $ws.Range("E51") "  +2.56%  "
